# Add a new "Run 50" results column right after "Run 49" (which currently
# sits in column AZ, immediately before the "Mean" summary column).
# This shifts the existing "Mean" column one position to the right, from
# AZ to BA, and fills column AZ with the new "Run 50" data and column BA
# with the recomputed "Mean" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 14

# 1. Copy the formatting of the current "Mean" column (AZ, including the
#    bold/centered/bordered header style) into the new last column (BA)
#    so the relocated "Mean" column keeps looking the same.
$ws.Range("AZ1:AZ$lastRow").Copy()
$ws.Range("BA1:BA$lastRow").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# 2. Headers: AZ1 becomes "Run 50" (reusing the column formerly used for
#    "Mean"), BA1 becomes the new "Mean" header.
$ws.Range("AZ1").Value = "Run 50"
$ws.Range("BA1").Value = "Mean"

# 3. Data rows: AZ gets the new "Run 50" results, BA gets the recomputed
#    "Mean" values (same figure for every row, as in the source data).
$run50Value = 139.54959404
$meanValue  = 131.67399209

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 52).Value = $run50Value
    $ws.Cells.Item($r, 53).Value = $meanValue
}
